$d = $word.ActiveDocument

# 1. Rename the "Enregistrement rotation d'une image" task.
$d.Content.Find.Execute(
    "Enregistrement rotation d’une image", $true, $false, $false, $false, $false,
    $true, 1, $false, "Améliorer les rotations d’images", 2)

# 2. Locate the "drag and drop" list item, then append two new paragraphs
#    after it: a new task, and a trailing empty (no-number) paragraph.
$count = $d.Paragraphs.Count
$dragParagraph = $null
for ($i = 1; $i -le $count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*drag and drop*") {
        $dragParagraph = $candidate
    }
}

$dragParagraph.Range.InsertParagraphAfter()
$navIndex = $dragParagraph.Index + 1
$navParagraph = $d.Paragraphs.Item($navIndex)
$navParagraph.Range.Text = "Ajouter le menu de navigation"

$navParagraph.Range.InsertParagraphAfter()
$trailingIndex = $navParagraph.Index + 1
$trailingParagraph = $d.Paragraphs.Item($trailingIndex)
$trailingParagraph.Range.ListFormat.RemoveNumbers()
$trailingParagraph.LeftIndent = 36
$trailingParagraph.FirstLineIndent = -0.01
